# Auto-generated edit script applying the Spriggan_Profits commit diff
# Updates currentAveragePrice / LevePrice / LeveProfit columns (H-N)
# across the ALC, ARM, BSM, CRP, CUL, GSM, LTW and WVR sheets.

$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 273.16666
$ws.Range("I33").Value = 273.16666
$ws.Range("K33").Value = 273.16666
$ws.Range("M33").Value = -44.16665999999998
$ws.Range("H41").Value = 1319.0714
$ws.Range("I41").Value = 2404.75
$ws.Range("J41").Value = 884.8
$ws.Range("K41").Value = 2404.75
$ws.Range("L41").Value = 884.8
$ws.Range("M41").Value = -1964.75
$ws.Range("N41").Value = -1764.8
$ws.Range("H100").Value = 4219.2
$ws.Range("I100").Value = 3765.3333
$ws.Range("K100").Value = 3765.3333
$ws.Range("M100").Value = -3224.3333
$ws.Range("H138").Value = 3182.4814
$ws.Range("J138").Value = 4140.4546
$ws.Range("L138").Value = 12421.3638
$ws.Range("N138").Value = -22701.3638

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1226274
$ws.Range("J2").Value = 2141.6667
$ws.Range("L2").Value = 2141.6667
$ws.Range("N2").Value = -2367.6667
$ws.Range("H32").Value = 4828.294
$ws.Range("I32").Value = 4828.294
$ws.Range("K32").Value = 4828.294
$ws.Range("M32").Value = -4541.294
$ws.Range("H61").Value = 90913490
$ws.Range("I61").Value = 100004390
$ws.Range("K61").Value = 100004390
$ws.Range("M61").Value = -100004178
$ws.Range("H97").Value = 883.4706
$ws.Range("I97").Value = 877.7143
$ws.Range("K97").Value = 877.7143
$ws.Range("M97").Value = -381.7143
$ws.Range("H102").Value = 11112565
$ws.Range("J102").Value = 2000
$ws.Range("L102").Value = 2000
$ws.Range("N102").Value = -5244
$ws.Range("H116").Value = 1226274
$ws.Range("J116").Value = 2141.6667
$ws.Range("L116").Value = 2141.6667
$ws.Range("N116").Value = -6729.6667
$ws.Range("H132").Value = 1590767.4
$ws.Range("I132").Value = 2002977.9
$ws.Range("K132").Value = 6008933.699999999
$ws.Range("M132").Value = -6006403.699999999
$ws.Range("H136").Value = 90913490
$ws.Range("I136").Value = 100004390
$ws.Range("K136").Value = 300013170
$ws.Range("M136").Value = -300010620

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1226274
$ws.Range("J3").Value = 2141.6667
$ws.Range("L3").Value = 2141.6667
$ws.Range("N3").Value = -2369.6667
$ws.Range("H86").Value = 3325.4285
$ws.Range("J86").Value = 2603.3333
$ws.Range("L86").Value = 2603.3333
$ws.Range("N86").Value = -4849.3333
$ws.Range("H89").Value = 3325.4285
$ws.Range("J89").Value = 2603.3333
$ws.Range("L89").Value = 13016.6665
$ws.Range("N89").Value = -24248.6665

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 361.33334
$ws.Range("J7").Value = 496.6
$ws.Range("L7").Value = 496.6
$ws.Range("N7").Value = -722.6
$ws.Range("H22").Value = 20439.4
$ws.Range("I22").Value = 25424.25
$ws.Range("K22").Value = 25424.25
$ws.Range("M22").Value = -25074.25
$ws.Range("H86").Value = 4805.25
$ws.Range("I86").Value = 4740.6665
$ws.Range("K86").Value = 4740.6665
$ws.Range("M86").Value = -3617.6665
$ws.Range("H89").Value = 4805.25
$ws.Range("I89").Value = 4740.6665
$ws.Range("K89").Value = 23703.3325
$ws.Range("M89").Value = -18087.3325
$ws.Range("H99").Value = 1488.1111
$ws.Range("I99").Value = 1299
$ws.Range("K99").Value = 1299
$ws.Range("M99").Value = 199
$ws.Range("H107").Value = 743886.0600000001
$ws.Range("I107").Value = 906294.7
$ws.Range("K107").Value = 906294.7
$ws.Range("M107").Value = -904374.7
$ws.Range("H122").Value = 2576.8
$ws.Range("I122").Value = 2435.389
$ws.Range("K122").Value = 7306.167
$ws.Range("M122").Value = -4856.167
$ws.Range("H126").Value = 1488.1111
$ws.Range("I126").Value = 1299
$ws.Range("K126").Value = 3897
$ws.Range("M126").Value = -1427
$ws.Range("H132").Value = 144447660
$ws.Range("I132").Value = 144447660
$ws.Range("K132").Value = 433342980
$ws.Range("M132").Value = -433340450

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 725.2
$ws.Range("I3").Value = 156.5
$ws.Range("J3").Value = 3000
$ws.Range("K3").Value = 469.5
$ws.Range("L3").Value = 9000
$ws.Range("M3").Value = -357.5
$ws.Range("N3").Value = -9224
$ws.Range("H14").Value = 378.22223
$ws.Range("I14").Value = 378.22223
$ws.Range("K14").Value = 1134.66669
$ws.Range("M14").Value = -961.66669
$ws.Range("H108").Value = 1942.1666
$ws.Range("I108").Value = 1846.091
$ws.Range("K108").Value = 5538.272999999999
$ws.Range("M108").Value = -2658.272999999999
$ws.Range("H129").Value = 1396.4615
$ws.Range("I129").Value = 794.375
$ws.Range("J129").Value = 2359.8
$ws.Range("K129").Value = 2383.125
$ws.Range("L129").Value = 7079.400000000001
$ws.Range("M129").Value = 2616.875
$ws.Range("N129").Value = -17079.4
$ws.Range("H131").Value = 4753.6665
$ws.Range("J131").Value = 4798.6
$ws.Range("L131").Value = 14395.8
$ws.Range("N131").Value = -24475.8

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3151
$ws.Range("I80").Value = 2866
$ws.Range("J80").Value = 4006
$ws.Range("K80").Value = 2866
$ws.Range("L80").Value = 4006
$ws.Range("M80").Value = -1868
$ws.Range("N80").Value = -6002
$ws.Range("H83").Value = 3151
$ws.Range("I83").Value = 2866
$ws.Range("J83").Value = 4006
$ws.Range("K83").Value = 14330
$ws.Range("L83").Value = 20030
$ws.Range("M83").Value = -9338
$ws.Range("N83").Value = -30014
$ws.Range("H98").Value = 95120.28999999999
$ws.Range("J98").Value = 95120.28999999999
$ws.Range("L98").Value = 95120.28999999999
$ws.Range("N98").Value = -101110.29

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("M17").ClearContents()

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 38228.332
$ws.Range("J41").Value = 38228.332
$ws.Range("L41").Value = 38228.332
$ws.Range("N41").Value = -39008.332
$ws.Range("H62").Value = 10999.667
$ws.Range("I62").Value = 10999.5
$ws.Range("J62").Value = 11000
$ws.Range("K62").Value = 10999.5
$ws.Range("L62").Value = 11000
$ws.Range("M62").Value = -10375.5
$ws.Range("N62").Value = -12248
$ws.Range("H65").Value = 10999.667
$ws.Range("I65").Value = 10999.5
$ws.Range("J65").Value = 11000
$ws.Range("K65").Value = 54997.5
$ws.Range("L65").Value = 55000
$ws.Range("M65").Value = -51877.5
$ws.Range("N65").Value = -61240
$ws.Range("H122").Value = 5725.231
$ws.Range("I122").Value = 6349.909
$ws.Range("K122").Value = 19049.727
$ws.Range("M122").Value = -16599.727
$ws.Range("H140").Value = 87347.5
$ws.Range("J140").Value = 87347.5
$ws.Range("L140").Value = 87347.5
$ws.Range("N140").Value = -97707.5
